# "disabled next button all together"
#
# Semantic change: the seven "select_one <list>" prompts on the "survey"
# sheet that implement menu/navigation screens are switched from the
# select_one widget to the custom "menu" widget (type column values
# "select_one X" -> "menu X"). A matching "menu" -> "string" schema-type
# mapping row is added to the prompt_types sheet, and the active
# worksheet/selection bookmarks are moved from the "choices" sheet back to
# the "survey" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. survey sheet: change "select_one X" type cells to "menu X"
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("A3").Value  = "menu menu"
$survey.Range("A8").Value  = "menu diagnostic_menu"
$survey.Range("A12").Value = "menu ballard_menu1"
$survey.Range("A24").Value = "menu procedures_menu"
$survey.Range("A37").Value = "menu calculators_menu"
$survey.Range("A55").Value = "menu nutritions_menu"
$survey.Range("A49").Value = "menu medications_menu"

# ---------------------------------------------------------------------
# 2. prompt_types sheet: add a "menu" -> "string" schema-type mapping
# ---------------------------------------------------------------------
$promptTypes = $wb.Worksheets.Item("prompt_types")
$promptTypes.Range("A3").Value = "menu"
$promptTypes.Range("B3").Value = "string"

# ---------------------------------------------------------------------
# 3. Move the active tab / selection from "choices" back to "survey"
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Application.Goto($choices.Range("B35"))

$promptTypes.Application.Goto($promptTypes.Range("A4"))

$survey.Activate()
$survey.Application.Goto($survey.Range("A50"))
